$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows (original rows 2-4) where sending cluster was "ECs".
# This shifts the remaining six rows (previously rows 5-10) up to rows 2-7,
# and the used range / dimension shrinks from A1:T10 to A1:T7 automatically.
$ws.Rows.Item(2).Resize(3).Delete() | Out-Null

# Update the numeric columns (G:T) for the six remaining data rows with the
# refreshed TPM-derived values.
$ws.Cells.Item(2,7).Value = 2.015377
$ws.Cells.Item(2,8).Value = 6.046131
$ws.Cells.Item(2,9).Value = 0.7554960962715589
$ws.Cells.Item(2,10).Value = 0.7554960962715588
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.1944653333333334
$ws.Cells.Item(2,14).Value = 0.583396
$ws.Cells.Item(2,15).Value = 0.04942840076761122
$ws.Cells.Item(2,16).Value = 0.04942840076761121
$ws.Cells.Item(2,17).Value = 0.3919209600973333
$ws.Cells.Item(2,18).Value = 3.527288640876
$ws.Cells.Item(2,19).Value = 0.0373429638248764
$ws.Cells.Item(2,20).Value = 0.03734296382487639

$ws.Cells.Item(3,7).Value = 2.015377
$ws.Cells.Item(3,8).Value = 6.046131
$ws.Cells.Item(3,9).Value = 0.7554960962715589
$ws.Cells.Item(3,10).Value = 0.7554960962715588
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.1693566666666667
$ws.Cells.Item(3,14).Value = 0.50807
$ws.Cells.Item(3,15).Value = 0.04304638286515546
$ws.Cells.Item(3,16).Value = 0.04304638286515546
$ws.Cells.Item(3,17).Value = 0.3413175307966667
$ws.Cells.Item(3,18).Value = 3.07185777717
$ws.Cells.Item(3,19).Value = 0.03252137421323587
$ws.Cells.Item(3,20).Value = 0.03252137421323587

$ws.Cells.Item(4,7).Value = 2.015377
$ws.Cells.Item(4,8).Value = 6.046131
$ws.Cells.Item(4,9).Value = 0.7554960962715589
$ws.Cells.Item(4,10).Value = 0.7554960962715588
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 3.570461333333333
$ws.Cells.Item(4,14).Value = 10.711384
$ws.Cells.Item(4,15).Value = 0.9075252163672334
$ws.Cells.Item(4,16).Value = 0.9075252163672333
$ws.Cells.Item(4,17).Value = 7.195825650589332
$ws.Cells.Item(4,18).Value = 64.76243085530399
$ws.Cells.Item(4,19).Value = 0.6856317582334467
$ws.Cells.Item(4,20).Value = 0.6856317582334465

$ws.Cells.Item(5,7).Value = 0.6522436666666667
$ws.Cells.Item(5,8).Value = 1.956731
$ws.Cells.Item(5,9).Value = 0.2445039037284412
$ws.Cells.Item(5,10).Value = 0.2445039037284411
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.1944653333333334
$ws.Cells.Item(5,14).Value = 0.583396
$ws.Cells.Item(5,15).Value = 0.04942840076761122
$ws.Cells.Item(5,16).Value = 0.04942840076761121
$ws.Cells.Item(5,17).Value = 0.1268387820528889
$ws.Cells.Item(5,18).Value = 1.141549038476
$ws.Cells.Item(5,19).Value = 0.01208543694273482
$ws.Cells.Item(5,20).Value = 0.01208543694273482

$ws.Cells.Item(6,7).Value = 0.6522436666666667
$ws.Cells.Item(6,8).Value = 1.956731
$ws.Cells.Item(6,9).Value = 0.2445039037284412
$ws.Cells.Item(6,10).Value = 0.2445039037284411
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.1693566666666667
$ws.Cells.Item(6,14).Value = 0.50807
$ws.Cells.Item(6,15).Value = 0.04304638286515546
$ws.Cells.Item(6,16).Value = 0.04304638286515546
$ws.Cells.Item(6,17).Value = 0.1104618132411111
$ws.Cells.Item(6,18).Value = 0.99415631917
$ws.Cells.Item(6,19).Value = 0.01052500865191959
$ws.Cells.Item(6,20).Value = 0.01052500865191959

$ws.Cells.Item(7,7).Value = 0.6522436666666667
$ws.Cells.Item(7,8).Value = 1.956731
$ws.Cells.Item(7,9).Value = 0.2445039037284412
$ws.Cells.Item(7,10).Value = 0.2445039037284411
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 3.570461333333333
$ws.Cells.Item(7,14).Value = 10.711384
$ws.Cells.Item(7,15).Value = 0.9075252163672334
$ws.Cells.Item(7,16).Value = 0.9075252163672333
$ws.Cells.Item(7,17).Value = 2.328810791744889
$ws.Cells.Item(7,18).Value = 20.959297125704
$ws.Cells.Item(7,19).Value = 0.2218934581337868
$ws.Cells.Item(7,20).Value = 0.2218934581337867
